{"js": "// Replace the inline \"road buffer / setback\" illustration picture with a\n// hyperlink run that points at the image's URL on ura.gov.sg, matching the\n// commit that swapped the embedded <w:drawing> picture for a\n// <w:hyperlink><w:r rStyle=\"Hyperlink\">...</w:r></w:hyperlink> run.\n\nconst url =\n  \"https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Others/HMC02_Road_Buffer_Setbacks.jpg?h=100%25&w=100%25\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Find the paragraph that holds the inline picture (the one right after the\n// \"Building Setback from Boundary\" Heading3, styled FirstParagraph).\nlet pictureParagraph = null;\nfor (const p of paragraphs.items) {\n  const pics = p.inlinePictures;\n  pics.load(\"items\");\n  await context.sync();\n  if (pics.items.length > 0) {\n    pictureParagraph = p;\n    break;\n  }\n}\n\nif (!pictureParagraph) {\n  throw new Error(\"Could not find the paragraph containing the inline picture\");\n}\n\nconst pics = pictureParagraph.inlinePictures;\npics.load(\"items\");\nawait context.sync();\nconst picture = pics.items[0];\n\n// Grab the picture's range so we can insert replacement text in the very\n// same spot, then delete the picture itself.\nconst pictureRange = picture.getRange();\npicture.delete();\nawait context.sync();\n\n// Insert the URL as plain text, then turn that text into a hyperlink \u2014\n// Word applies the built-in \"Hyperlink\" character style automatically.\nconst textRange = pictureRange.insertText(url, \"Replace\");\nawait context.sync();\n\ntextRange.hyperlink = url;\nawait context.sync();\n", "ps1": "# Replace the inline \"road buffer / setback\" illustration picture with a\n# hyperlink run that points at the image's URL on ura.gov.sg, matching the\n# commit that swapped the embedded picture for a\n# <w:hyperlink><w:r rStyle=\"Hyperlink\">...</w:r></w:hyperlink> run.\n\n$d = $word.ActiveDocument\n$url = \"https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Others/HMC02_Road_Buffer_Setbacks.jpg?h=100%25&w=100%25\"\n\n$pic = $d.InlineShapes.Item(1)\n$picRange = $pic.Range\n$pic.Delete()\n\n$d.Hyperlinks.Add($picRange, $url, $null, $null, $url)\n"}
